$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) values for rows 2-11
$ws.Range("C2").Value = -1.054401637039784
$ws.Range("D2").Value = 0.3031419272028697

$ws.Range("C3").Value = -0.5481820004710556
$ws.Range("D3").Value = 0.5890867629772574

$ws.Range("C4").Value = -0.3229303236714468
$ws.Range("D4").Value = 0.7497972314666175

$ws.Range("C5").Value = -1.652965517199736
$ws.Range("D5").Value = 0.1125389990266088
$ws.Range("G5").Value = "No"

$ws.Range("C6").Value = 0.4925252054464042
$ws.Range("D6").Value = 0.6272276353813044

$ws.Range("C7").Value = 0.4499715141939513
$ws.Range("D7").Value = 0.6571334895511902

$ws.Range("C8").Value = -0.7319989567584946
$ws.Range("D8").Value = 0.471895771698865

$ws.Range("C9").Value = 0.05939148799542811
$ws.Range("D9").Value = 0.953176491385938

$ws.Range("C10").Value = -1.169870763767922
$ws.Range("D10").Value = 0.2545712277780448

$ws.Range("C11").Value = -1.239461460582337
$ws.Range("D11").Value = 0.2282344319631306

$wb.Save()
